$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Split the big intro paragraph into four paragraphs:
#      "Make the normal spell to grow..."
#      "Create a mask for the vision when teleporting."
#      (blank paragraph)
#      "You have to create a particle effect... press the grip."
# ---------------------------------------------------------------------------
$old1 = "You have to create a particle effect for the platform that is the target of the teleportation. Make the particle effect disappear if the pointer leaves the platform. Make the seal active if you press the grip. Create a mask for the vision when teleporting."
$new1 = "Make the normal spell to grow with the pressed trigger and release it when you release the trigger.^pCreate a mask for the vision when teleporting.^p^pYou have to create a particle effect for the platform that is the target of the teleportation. Make the particle effect disappear if the pointer leaves the platform. Make the seal active if you press the grip."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Prepend the "11/26/2016 " date stamp as its own run ahead of the moved sentence.
$r = $d.Content
$r.Find.Execute("You have to create a particle effect", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($r.Start, $r.Start)
$insPoint.InsertBefore("11/26/2016 ")

# ---------------------------------------------------------------------------
# 2. Small wording/typo fixes further down in the log.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The came is still working.", $false, $false, $false, $false, $false, $true, 1, $false, "The game is still working.", 2) | Out-Null
$d.Content.Find.Execute("The hand are uploaded at runtime.", $false, $false, $false, $false, $false, $true, 1, $false, "The hands are uploaded at runtime.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Relocate the hidden "_GoBack" bookmark (Word re-stamps it at the last
#    edited spot automatically; put it explicitly where the diff wants it -
#    right after "that will send it to" in the 11/11/2016 paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$anchor = $d.Content
$anchor.Find.Execute("that will send it to", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
